$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 356.4
$ws.Cells.Item(2, 9).Value = 145.5
$ws.Cells.Item(2, 11).Value = 145.5
$ws.Cells.Item(2, 13).Value = -32.5

$ws.Cells.Item(15, 8).Value = 1768.1708
$ws.Cells.Item(15, 9).Value = 1768.1708
$ws.Cells.Item(15, 11).Value = 5304.512400000001
$ws.Cells.Item(15, 13).Value = -5135.512400000001

$ws.Cells.Item(51, 8).Value = 2997.5

$ws.Cells.Item(58, 8).Value = 1235
$ws.Cells.Item(58, 10).Value = 1770.909
$ws.Cells.Item(58, 12).Value = 5312.727000000001
$ws.Cells.Item(58, 14).Value = -5612.727000000001

$ws.Cells.Item(61, 8).Value = 14228
$ws.Cells.Item(61, 10).Value = 1000
$ws.Cells.Item(61, 12).Value = 3000
$ws.Cells.Item(61, 14).Value = -3344

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).Value = $null

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).Value = $null

$ws.Cells.Item(80, 8).Value = 637.75
$ws.Cells.Item(80, 9).Value = 670.1667
$ws.Cells.Item(80, 10).Value = 605.3333
$ws.Cells.Item(80, 11).Value = 2010.5001
$ws.Cells.Item(80, 12).Value = 1815.9999
$ws.Cells.Item(80, 13).Value = -1012.5001
$ws.Cells.Item(80, 14).Value = -3811.9999

$ws.Cells.Item(83, 8).Value = 637.75
$ws.Cells.Item(83, 9).Value = 670.1667
$ws.Cells.Item(83, 10).Value = 605.3333
$ws.Cells.Item(83, 11).Value = 6031.5003
$ws.Cells.Item(83, 12).Value = 5447.9997
$ws.Cells.Item(83, 13).Value = -1039.5003
$ws.Cells.Item(83, 14).Value = -15431.9997

$ws.Cells.Item(88, 8).Value = 2389.5
$ws.Cells.Item(88, 9).Value = 2209.5
$ws.Cells.Item(88, 10).Value = 2479.5
$ws.Cells.Item(88, 11).Value = 2209.5
$ws.Cells.Item(88, 12).Value = 2479.5
$ws.Cells.Item(88, 13).Value = -1803.5
$ws.Cells.Item(88, 14).Value = -3291.5

$ws.Cells.Item(91, 8).Value = 2389.5
$ws.Cells.Item(91, 9).Value = 2209.5
$ws.Cells.Item(91, 10).Value = 2479.5
$ws.Cells.Item(91, 11).Value = 2209.5
$ws.Cells.Item(91, 12).Value = 2479.5
$ws.Cells.Item(91, 13).Value = -805.5
$ws.Cells.Item(91, 14).Value = -5287.5

$ws.Cells.Item(100, 8).Value = 1249.6666
$ws.Cells.Item(100, 9).Value = 874.5
$ws.Cells.Item(100, 11).Value = 874.5
$ws.Cells.Item(100, 13).Value = -333.5

$ws.Cells.Item(138, 8).Value = 2357.5
$ws.Cells.Item(138, 9).Value = 1180
$ws.Cells.Item(138, 11).Value = 3540
$ws.Cells.Item(138, 13).Value = 1600

$ws.Cells.Item(141, 8).Value = 1122
$ws.Cells.Item(141, 9).Value = 1122
$ws.Cells.Item(141, 11).Value = 3366
$ws.Cells.Item(141, 13).Value = 1814


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).Value = $null

$ws.Cells.Item(32, 8).Value = 3336018
$ws.Cells.Item(32, 9).Value = 621.72
$ws.Cells.Item(32, 11).Value = 621.72
$ws.Cells.Item(32, 13).Value = -334.72

$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).Value = $null

$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).Value = $null

$ws.Cells.Item(132, 8).Value = 4062.9167
$ws.Cells.Item(132, 9).Value = 4069.5454
$ws.Cells.Item(132, 11).Value = 12208.6362
$ws.Cells.Item(132, 13).Value = -9678.636200000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 10954099
$ws.Cells.Item(7, 10).Value = 7187507
$ws.Cells.Item(7, 12).Value = 7187507
$ws.Cells.Item(7, 14).Value = -7187733

$ws.Cells.Item(8, 8).Value = 1118
$ws.Cells.Item(8, 9).Value = 1118
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 1118
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -978
$ws.Cells.Item(8, 14).Value = $null

$ws.Cells.Item(11, 8).Value = 866
$ws.Cells.Item(11, 9).Value = 1433.3334
$ws.Cells.Item(11, 10).Value = 15
$ws.Cells.Item(11, 11).Value = 1433.3334
$ws.Cells.Item(11, 12).Value = 15
$ws.Cells.Item(11, 13).Value = -1293.3334
$ws.Cells.Item(11, 14).Value = -295

$ws.Cells.Item(12, 8).Value = 263.33334
$ws.Cells.Item(12, 9).Value = 337.5
$ws.Cells.Item(12, 10).Value = 115
$ws.Cells.Item(12, 11).Value = 337.5
$ws.Cells.Item(12, 12).Value = 115
$ws.Cells.Item(12, 13).Value = -169.5
$ws.Cells.Item(12, 14).Value = -451

$ws.Cells.Item(16, 8).Value = 438.66666
$ws.Cells.Item(16, 9).Value = 438.66666
$ws.Cells.Item(16, 11).Value = 438.66666
$ws.Cells.Item(16, 13).Value = -268.66666

$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = $null

$ws.Cells.Item(105, 8).Value = 5348811
$ws.Cells.Item(105, 9).Value = 7576765.5
$ws.Cells.Item(105, 10).Value = 1720
$ws.Cells.Item(105, 11).Value = 7576765.5
$ws.Cells.Item(105, 12).Value = 1720
$ws.Cells.Item(105, 13).Value = -7575018.5
$ws.Cells.Item(105, 14).Value = -5214

$ws.Cells.Item(134, 8).Value = 4689.905
$ws.Cells.Item(134, 9).Value = 1511.8125
$ws.Cells.Item(134, 11).Value = 4535.4375
$ws.Cells.Item(134, 13).Value = -2000.4375


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3200.068
$ws.Cells.Item(31, 9).Value = 1895.8572
$ws.Cells.Item(31, 10).Value = 4390.8696
$ws.Cells.Item(31, 11).Value = 1895.8572
$ws.Cells.Item(31, 12).Value = 4390.8696
$ws.Cells.Item(31, 13).Value = -1600.8572
$ws.Cells.Item(31, 14).Value = -4980.8696

$ws.Cells.Item(32, 8).Value = 1079.2
$ws.Cells.Item(32, 9).Value = 921.3333
$ws.Cells.Item(32, 11).Value = 921.3333
$ws.Cells.Item(32, 13).Value = -605.3333

$ws.Cells.Item(34, 8).Value = 3200.068
$ws.Cells.Item(34, 9).Value = 1895.8572
$ws.Cells.Item(34, 10).Value = 4390.8696
$ws.Cells.Item(34, 11).Value = 1895.8572
$ws.Cells.Item(34, 12).Value = 4390.8696
$ws.Cells.Item(34, 13).Value = -1693.8572
$ws.Cells.Item(34, 14).Value = -4794.8696

$ws.Cells.Item(45, 8).Value = 39999
$ws.Cells.Item(45, 10).Value = 39999
$ws.Cells.Item(45, 12).Value = 39999
$ws.Cells.Item(45, 14).Value = -41185

$ws.Cells.Item(103, 8).Value = 16246.75
$ws.Cells.Item(103, 9).Value = 16246.75
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 16246.75
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -15074.75
$ws.Cells.Item(103, 14).Value = $null

$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = $null
$ws.Cells.Item(105, 14).Value = $null

$ws.Cells.Item(106, 8).Value = 105000
$ws.Cells.Item(106, 10).Value = 105000
$ws.Cells.Item(106, 12).Value = 105000
$ws.Cells.Item(106, 14).Value = -107524

$ws.Cells.Item(132, 8).Value = 2228.8845
$ws.Cells.Item(132, 9).Value = 2223.9565
$ws.Cells.Item(132, 11).Value = 6671.869499999999
$ws.Cells.Item(132, 13).Value = -4141.869499999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 259.8
$ws.Cells.Item(38, 9).Value = 291.375
$ws.Cells.Item(38, 11).Value = 874.125
$ws.Cells.Item(38, 13).Value = -527.125

$ws.Cells.Item(107, 8).Value = 490
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 490
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 1470
$ws.Cells.Item(107, 13).Value = $null
$ws.Cells.Item(107, 14).Value = -5310


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 4994
$ws.Cells.Item(57, 9).Value = 4994
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 4994
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = -4174
$ws.Cells.Item(57, 14).Value = $null

$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 13).Value = $null

$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 13).Value = $null

$ws.Cells.Item(102, 8).Value = 2249.7896
$ws.Cells.Item(102, 9).Value = 2097
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 2097
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -475
$ws.Cells.Item(102, 14).Value = -8244

$ws.Cells.Item(122, 8).Value = 1761
$ws.Cells.Item(122, 9).Value = 1147.5
$ws.Cells.Item(122, 10).Value = 2742.6
$ws.Cells.Item(122, 11).Value = 3442.5
$ws.Cells.Item(122, 12).Value = 8227.799999999999
$ws.Cells.Item(122, 13).Value = -992.5
$ws.Cells.Item(122, 14).Value = -13127.8

$ws.Cells.Item(126, 8).Value = 2803.4546
$ws.Cells.Item(126, 9).Value = 2784.8
$ws.Cells.Item(126, 10).Value = 2990
$ws.Cells.Item(126, 11).Value = 8354.400000000001
$ws.Cells.Item(126, 12).Value = 8970
$ws.Cells.Item(126, 13).Value = -5884.400000000001
$ws.Cells.Item(126, 14).Value = -13910

$ws.Cells.Item(132, 8).Value = 26749
$ws.Cells.Item(132, 9).Value = 31267.686
$ws.Cells.Item(132, 11).Value = 93803.058
$ws.Cells.Item(132, 13).Value = -91273.058


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2800
$ws.Cells.Item(7, 10).Value = 2800
$ws.Cells.Item(7, 12).Value = 2800
$ws.Cells.Item(7, 14).Value = -3024

$ws.Cells.Item(53, 8).Value = 5999.5
$ws.Cells.Item(53, 9).Value = 5999.5
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 5999.5
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = -5481.5
$ws.Cells.Item(53, 14).Value = $null

$ws.Cells.Item(61, 8).Value = 66670460
$ws.Cells.Item(61, 9).Value = 125001240
$ws.Cells.Item(61, 11).Value = 125001240
$ws.Cells.Item(61, 13).Value = -125001038

$ws.Cells.Item(113, 8).Value = 66670460
$ws.Cells.Item(113, 9).Value = 125001240
$ws.Cells.Item(113, 11).Value = 125001240
$ws.Cells.Item(113, 13).Value = -124999070

$ws.Cells.Item(126, 8).Value = 2800
$ws.Cells.Item(126, 10).Value = 2800
$ws.Cells.Item(126, 12).Value = 8400
$ws.Cells.Item(126, 14).Value = -13340

$ws.Cells.Item(136, 8).Value = 2107.4167
$ws.Cells.Item(136, 9).Value = 1911.1111
$ws.Cells.Item(136, 10).Value = 2696.3333
$ws.Cells.Item(136, 11).Value = 5733.3333
$ws.Cells.Item(136, 12).Value = 8088.999899999999
$ws.Cells.Item(136, 13).Value = -3183.3333
$ws.Cells.Item(136, 14).Value = -13188.9999

